$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells to preserve original text representation
$textCells = @('D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D12', 'D13', 'D15', 'D17', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.161.74'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '2.904.56'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '347.55'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = '106.57'
$ws.Range('E6').Value = '  -5.99%  '
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  -2.42%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.603'
$ws.Range('E9').Value = '  -2.86%  '
$ws.Range('D10').Value = '37.35'
$ws.Range('E10').Value = '  -5.48%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = '0.0844'
$ws.Range('E12').Value = '  -3.85%  '
$ws.Range('D13').Value = '18.75'
$ws.Range('E13').Value = '  -6.46%  '
$ws.Range('D14').Value = '3.363.17'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').Value = '7.52'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').Value = '2.972.73'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = '0.951'
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '51.090.57'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').Value = '3.38'
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').Value = '7.34'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('D21').Value = '13.27'
$ws.Range('E21').Value = '  -6.46%  '
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('D23').Value = '68.50'
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').Value = '259.89'
$ws.Range('E24').Value = '  -3.34%  '
$ws.Range('D25').Value = '2.67'
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('E26').Value = '  -5.04%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '26.15'
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('D29').Value = '7.37'
$ws.Range('E29').Value = '  +5.35%  '
$ws.Range('D30').Value = '0.104'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').Value = '10.12'
$ws.Range('E31').Value = '  -4.79%  '
$ws.Range('D32').Value = '6.05'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = '35.19'
$ws.Range('E33').Value = '  -5.20%  '
$ws.Range('D34').Value = '2.11'
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('D35').Value = '50.14'
$ws.Range('E35').Value = '  -5.56%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = '0.0422'
$ws.Range('E37').Value = '  -6.96%  '
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').Value = '  -8.19%  '
$ws.Range('D39').Value = '17.49'
$ws.Range('E39').Value = '  -6.18%  '
$ws.Range('D40').Value = '1.92'
$ws.Range('E40').Value = '  -6.10%  '
$ws.Range('D41').Value = '2.61'
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').Value = '22.12'
$ws.Range('E43').Value = '  -3.41%  '
$ws.Range('D44').Value = '119.23'
$ws.Range('E44').Value = '  +6.81%  '
$ws.Range('E45').Value = '  -2.85%  '
$ws.Range('D46').Value = '2.082.22'
$ws.Range('E46').Value = '  -4.97%  '
$ws.Range('D47').Value = '3.27'
$ws.Range('E47').Value = '  -6.45%  '
$ws.Range('D48').Value = '2.25'
$ws.Range('E48').Value = '  -10.63%  '
$ws.Range('D49').Value = '0.238'
$ws.Range('E49').Value = '  -4.45%  '
$ws.Range('D50').Value = '0.0334'
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('B51').Value = 'SEI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D51').Value = '0.883'
$ws.Range('E51').Value = '  -7.03%  '
